# "updated myodds new leagues"
# Refresh the simulated match inputs (columns A/B) on the score_newleagues
# sheet. Columns C:F already contain (or, for newly-populated rows, will
# receive) the ROUND / CONCATENATE / VLOOKUP formulas that derive the
# correctscore grid from A/B, so Excel recalculates them automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("score_newleagues")

# New random-odds inputs for rows 1-20 and 22-30. Row 21 is a special
# case (a literal #N/A error pair) and is handled separately below.
$rows = @(
    @{ Row = 1; A = "2.0599788175540001"; B = "1.5264174314880001" },
    @{ Row = 2; A = "1.7114861604490001"; B = "0.96461101573199992" },
    @{ Row = 3; A = "0.51698056208400012"; B = "0.89010434368800007" },
    @{ Row = 4; A = "0.98688561606599989"; B = "1.7802086873760001" },
    @{ Row = 5; A = "0.64134129805199991"; B = "1.5216906450899998" },
    @{ Row = 6; A = "1.4637469359799999"; B = "1.3737490213680001" },
    @{ Row = 7; A = "0.959790776304"; B = "1.510991259516" },
    @{ Row = 8; A = "2.2513192303710001"; B = "1.3737934663200002" },
    @{ Row = 9; A = "0.48245252276400002"; B = "0.70089542107200009" },
    @{ Row = 10; A = "0.62656977042000006"; B = "0.20226280620000001" },
    @{ Row = 11; A = "0.82703289779400002"; B = "1.1571516472319998" },
    @{ Row = 12; A = "1.9249130937119998"; B = "1.2062852598640001" },
    @{ Row = 13; A = "1.9249130937119998"; B = "1.85703242322" },
    @{ Row = 14; A = "2.9459542338900002"; B = "0.89496863024999995" },
    @{ Row = 15; A = "0.92586149883000002"; B = "0.75017223417500001" },
    @{ Row = 16; A = "1.6272850838399999"; B = "1.2146574247649999" },
    @{ Row = 17; A = "0.71096279747600011"; B = "1.0303117660260002" },
    @{ Row = 18; A = "2.7810414022230003"; B = "0.76320871574400007" },
    @{ Row = 19; A = "1.0067565649699999"; B = "1.780820336736" },
    @{ Row = 20; A = "1.7230218683160001"; B = "0.89010434368800007" },
    @{ Row = 22; A = "3.5669301409279996"; B = "0.85704829972800012" },
    @{ Row = 23; A = "1.6559624088"; B = "0.85704829972800012" },
    @{ Row = 24; A = "0.93432694555199991"; B = "3.0474013594560003" },
    @{ Row = 25; A = "0.67779683625599996"; B = "0.74582254012500016" },
    @{ Row = 26; A = "2.1087672397919999"; B = "1.403944597575" },
    @{ Row = 27; A = "2.6717770378050001"; B = "1.0176116209920001" },
    @{ Row = 28; A = "2.0935957716360001"; B = "0.44500217160000005" },
    @{ Row = 29; A = "1.0592530249949998"; B = "0.97139950687999999" },
    @{ Row = 30; A = "0.73647826441199993"; B = "1.250357644795" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("A$r").Value = [double]$item.A
    $ws.Range("B$r").Value = [double]$item.B
}

# Row 21 keeps a literal #N/A error pair in A21/B21 (not a formula) -
# mirrors a pasted VLOOKUP failure elsewhere in the source data.
$ws.Range("A21").Value = "#N/A"
$ws.Range("B21").Value = "#N/A"

# Rows 14-30 previously had no formulas at all in C:F (they were blank
# placeholder rows). Populate them with the same formulas used by the
# rows above so the correct-score grid extends down to row 30.
for ($r = 14; $r -le 30; $r++) {
    $ws.Range("C$r").Formula = "=ROUND(A$r,0)"
    $ws.Range("D$r").Formula = "=ROUND(B$r,0)"
    $ws.Range("E$r").Formula = '=CONCATENATE(C' + $r + ',"-",D' + $r + ')'
    $ws.Range("F$r").Formula = '=VLOOKUP(E' + $r + ',cs_lookupnewleagues!$A$2:$B$54,2,FALSE)'
}

# Restore the view: scrolled to show row 4 at top, with E1:E30 selected.
$ws.Activate()
$ws.Range("E1:E30").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
